$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (D) / Volume(1h) (E) values per row, taken from the updated crypto feed.
$updates = @{
    2 = @{ D = "29.259.16"; E = "  -0.56%  " }
    3 = @{ D = "1.863.16"; E = "  -0.75%  " }
    4 = @{ E = "  +0.01%  " }
    5 = @{ D = "0.7144"; E = "  -0.56%  " }
    6 = @{ D = "240.99"; E = "  +0.34%  " }
    7 = @{ E = "  +0.08%  " }
    8 = @{ D = "0.07743"; E = "  -1.13%  " }
    9 = @{ D = "0.3085"; E = "  -0.92%  " }
    10 = @{ D = "24.97"; E = "  +0.02%  " }
    11 = @{ D = "0.08330"; E = "  +0.99%  " }
    12 = @{ D = "1.882.17"; E = "  +0.07%  " }
    13 = @{ D = "0.7167"; E = "  -1.45%  " }
    14 = @{ D = "5.214"; E = "  -1.38%  " }
    15 = @{ E = "  -0.36%  " }
    16 = @{ D = "29.270.48"; E = "  -0.79%  " }
    17 = @{ D = "5.972"; E = "  +0.76%  " }
    18 = @{ D = "243.10"; E = "  -0.91%  " }
    19 = @{ D = "0.000007833"; E = "  -0.61%  " }
    20 = @{ D = "2.132.12"; E = "  -0.47%  " }
    21 = @{ D = "13.17"; E = "  -0.90%  " }
    22 = @{ D = "1.000"; E = "  +0.05%  " }
    23 = @{ D = "7.920"; E = "  -0.59%  " }
    24 = @{ D = "1.001"; E = "  +0.10%  " }
    25 = @{ D = "0.1598"; E = "  +1.06%  " }
    26 = @{ D = "163.30"; E = "  -0.34%  " }
    27 = @{ D = "8.899"; E = "  -1.61%  " }
    28 = @{ E = "  +1.57%  " }
    29 = @{ D = "1.344"; E = "  -1.49%  " }
    30 = @{ D = "1.501"; E = "  +1.17%  " }
    31 = @{ D = "4.429"; E = "  +1.05%  " }
    32 = @{ D = "4.267"; E = "  +2.77%  " }
    33 = @{ D = "0.05170"; E = "  -2.09%  " }
    34 = @{ D = "0.8318"; E = "  +15.28%  " }
    35 = @{ E = "  -0.51%  " }
    36 = @{ D = "1.173"; E = "  -2.24%  " }
    37 = @{ D = "2.685"; E = "  +0.23%  " }
    38 = @{ D = "0.01857"; E = "  -0.50%  " }
    39 = @{ D = "2.693"; E = "  -1.20%  " }
    40 = @{ D = "1.165.44"; E = "  -5.97%  " }
    41 = @{ D = "6.199"; E = "  +1.63%  " }
    42 = @{ D = "0.8958"; E = "  -1.00%  " }
    43 = @{ D = "72.83"; E = "  -1.13%  " }
    44 = @{ E = "  +0.04%  " }
    45 = @{ D = "101.73"; E = "  -1.59%  " }
    46 = @{ D = "2.029.69"; E = "  +0.30%  " }
    47 = @{ D = "0.5180"; E = "  -2.95%  " }
    48 = @{ D = "1.786"; E = "  +1.42%  " }
    49 = @{ D = "9.363"; E = "  +0.97%  " }
    50 = @{ D = "7.078"; E = "  +0.03%  " }
    51 = @{ D = "0.4283"; E = "  -1.09%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        # Price text can look numeric (e.g. "1.000", "0.7144"); force Text format so
        # Excel does not silently coerce it to a number, then drop back to the default
        # (un-styled) look once the literal text is safely stored.
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $vals["D"]
        $ws.Range("D$row").Style = "Normal"
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
